# Update the "想去人数" (interested-count) column F values on the
# regenerated data sheets, matching the newly published gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1252
$ws1.Range("F5").Value = 5083
$ws1.Range("F7").Value = 9705
$ws1.Range("F11").Value = 28

# Sheet "演出" (shows)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 13

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1252
$ws4.Range("F7").Value = 5083
$ws4.Range("F9").Value = 13
$ws4.Range("F10").Value = 9705
$ws4.Range("F16").Value = 28
